# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted at row 150 (pushing the
# existing rows 150-157 down to 151-158), and the new row is populated
# with its data (Perú origin, 2021-11-09 / serial 44509, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 150; existing rows 150-157 move to 151-158.
$ws.Rows.Item(150).Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Cells.Item(150, 1).Value = 9
$ws.Cells.Item(150, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(150, 3).Value = "Metropolitana"
$ws.Cells.Item(150, 4).Value = 44509
$ws.Cells.Item(150, 5).Value = 13
$ws.Cells.Item(150, 6).Value = 100112030
$ws.Cells.Item(150, 7).Value = "Poroto granado"
$ws.Cells.Item(150, 8).Value = "Sin especificar"
$ws.Cells.Item(150, 9).Value = "Primera"
$ws.Cells.Item(150, 10).Value = 25
$ws.Cells.Item(150, 11).Value = 34000
$ws.Cells.Item(150, 12).Value = 36000
$ws.Cells.Item(150, 13).Value = 34960
$ws.Cells.Item(150, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(150, 15).Value = "Perú"
$ws.Cells.Item(150, 16).Value = 1398
$ws.Cells.Item(150, 17).Value = 25
$ws.Cells.Item(150, 18).Value = "Hortaliza"
